$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text by pre-formatting the
# data range (rows 2-51, columns D:E) as Text before assigning values,
# then restore the default (unstyled) look so no stray number formats remain.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.684.77'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '2.121.46'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('D4').Value = '1.012'
$ws.Range('D5').Value = '337.86'
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').Value = '0.5262'
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('D8').Value = '0.4560'
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('D9').Value = '54.98'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').Value = '0.09131'
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '2.119.31'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = '6.865'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').Value = '8.170'
$ws.Range('E15').Value = '  +5.88%  '
$ws.Range('D16').Value = '0.00001177'
$ws.Range('E16').Value = '  +4.61%  '
$ws.Range('D17').Value = '97.28'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '1.013'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('D20').Value = '19.51'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').Value = '6.331'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '30.755.00'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').Value = '12.92'
$ws.Range('E24').Value = '  +5.10%  '
$ws.Range('D25').Value = '2.367'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').Value = '2.363.94'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').Value = '22.44'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').Value = '165.46'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('D29').Value = '2.574'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '134.74'
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').Value = '1.214'
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').Value = '1.666'
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').Value = '6.386'
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('D35').Value = '3.944'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('D36').Value = '10.64'
$ws.Range('E36').Value = '  +1.63%  '
$ws.Range('D37').Value = '5.898'
$ws.Range('E37').Value = '  +7.77%  '
$ws.Range('D38').Value = '0.02655'
$ws.Range('E38').Value = '  +3.31%  '
$ws.Range('D39').Value = '0.06894'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').Value = '0.2329'
$ws.Range('D41').Value = '12.72'
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').Value = '0.6932'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '15.38'
$ws.Range('E43').Value = '  +9.54%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '1.261'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').Value = '0.6509'
$ws.Range('D46').Value = '2.329'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('E47').Value = '  +22.03%  '
$ws.Range('D48').Value = '3.699'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D50').Value = '83.63'
$ws.Range('D51').Value = '0.07317'
$ws.Range('E51').Value = '  +3.85%  '

# Restore the original (default/no explicit number format) style so the
# cells match the source workbook formatting.
$ws.Range("D2:E51").Style = "Normal"

